# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# F2: 9264 -> 9300
# F4: 3    -> 13
# F5: 495  -> 500

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9300
    $ws.Range("F4").Value = 13
    $ws.Range("F5").Value = 500
}
